# feat: add 2022-Q1 data
#
# Before: 4 sheets -> 2021-Q2, 2021-Q3, 2021-Q4, 总计 (totals)
# After:  5 sheets -> 2021-Q2, 2021-Q3, 2021-Q4, 2022-Q1, 总计 (totals)
#
# The new "2022-Q1" sheet holds the same per-fund holdings table as the other
# quarter sheets, and the "总计" (totals) sheet gains a new leading row
# summarising 2022-Q1.
#
# To land on the same sheetId / r:id layout as the authors (2022-Q1 reuses
# the old "总计" sheet's slot, and a fresh sheet is appended at the very end
# named "总计"), we:
#   1. Rename the existing "总计" sheet to "2022-Q1" and overwrite its data.
#   2. Append a brand new sheet named "总计" at the end with the refreshed
#      totals table.
#   3. Re-activate the sheet that was originally active (2021-Q2) since
#      adding/renaming sheets shifts Excel's active-tab selection.

$wb = $excel.ActiveWorkbook

$firstSheet = $wb.Worksheets.Item(1)
$q4 = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# Step 1: turn the existing "总计" sheet into the new "2022-Q1" sheet
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# Headers (row 1). B1:D1 already carry the bold/boxed header style from the
# old "总计" sheet; copy that same style onto the new E1:H1 header cells.
$q1.Cells.Item(1,2).Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)

$q1.Cells.Item(1,2).Value = "基金代码"
$q1.Cells.Item(1,3).Value = "基金名称"
$q1.Cells.Item(1,4).Value = "基金规模"
$q1.Cells.Item(1,5).Value = "股票总仓位"
$q1.Cells.Item(1,6).Value = "仓位占比"
$q1.Cells.Item(1,7).Value = "持有市值(亿元)"
$q1.Cells.Item(1,8).Value = "仓位排名"

# Row 2: 501303 / A share class
$q1.Cells.Item(2,1).Value = 0
$q1.Cells.Item(2,2).Value = "'501303"
$q1.Cells.Item(2,2).ClearFormats()
$q1.Cells.Item(2,3).Value = "广发港股通恒生综合中型股指数(LOF)A"
$q1.Cells.Item(2,4).Value = "'0.34"
$q1.Cells.Item(2,4).ClearFormats()
$q1.Cells.Item(2,5).Value = "'92.39"
$q1.Cells.Item(2,5).ClearFormats()
$q1.Cells.Item(2,6).Value = "'0.94"
$q1.Cells.Item(2,6).ClearFormats()
$q1.Cells.Item(2,7).Value = "'0.0032"
$q1.Cells.Item(2,7).ClearFormats()
$q1.Cells.Item(2,8).Value = 10

# Row 3: 004996 / C share class
$q1.Cells.Item(3,1).Value = 1
$q1.Cells.Item(3,2).Value = "'004996"
$q1.Cells.Item(3,2).ClearFormats()
$q1.Cells.Item(3,3).Value = "广发港股通恒生综合中型股指数(LOF)C"
$q1.Cells.Item(3,4).Value = "'0.11"
$q1.Cells.Item(3,4).ClearFormats()
$q1.Cells.Item(3,5).Value = "'92.39"
$q1.Cells.Item(3,5).ClearFormats()
$q1.Cells.Item(3,6).Value = "'0.94"
$q1.Cells.Item(3,6).ClearFormats()
$q1.Cells.Item(3,7).Value = "'0.0010"
$q1.Cells.Item(3,7).ClearFormats()
$q1.Cells.Item(3,8).Value = 10

# The old "总计" sheet had a 4th row (2021-Q2 totals) that doesn't belong in
# the new holdings table - remove it entirely so the sheet is only 3 rows.
$q1.Rows.Item(4).Delete()

# ---------------------------------------------------------------------
# Step 2: append a brand-new "总计" sheet at the end with the updated
# totals table (2022-Q1 row added on top of the previous quarters).
# ---------------------------------------------------------------------
$lastIndex = $wb.Worksheets.Count
$total = $wb.Worksheets.Add($null, $wb.Worksheets.Item($lastIndex))
$total.Name = "总计"

# Match the page margins used by the other (pre-existing) sheets.
$total.PageSetup.LeftMargin = 54
$total.PageSetup.RightMargin = 54
$total.PageSetup.TopMargin = 72
$total.PageSetup.BottomMargin = 72
$total.PageSetup.HeaderMargin = 36
$total.PageSetup.FooterMargin = 36

# Match the header / index-column styling used on the other sheets.
$q4.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$q4.Cells.Item(2,1).Copy()
$total.Range("A2:A5").PasteSpecial(-4122)

$total.Cells.Item(1,2).Value = "日期"
$total.Cells.Item(1,3).Value = "持有数量(只)"
$total.Cells.Item(1,4).Value = "持有市值(亿元)"

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q1"
$total.Cells.Item(2,3).Value = 2
$total.Cells.Item(2,4).Value = 0

$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(3,2).Value = "2021-Q4"
$total.Cells.Item(3,3).Value = 2
$total.Cells.Item(3,4).Value = 0

$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(4,2).Value = "2021-Q3"
$total.Cells.Item(4,3).Value = 2
$total.Cells.Item(4,4).Value = 0

$total.Cells.Item(5,1).Value = 3
$total.Cells.Item(5,2).Value = "2021-Q2"
$total.Cells.Item(5,3).Value = 2
$total.Cells.Item(5,4).Value = 0

# ---------------------------------------------------------------------
# Step 3: restore the original active sheet/tab selection.
# ---------------------------------------------------------------------
$firstSheet.Activate()
